$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.483.34'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '3.396.71'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.55%  '
$ws.Range('D8').Value = '3.399.27'
$ws.Range('E8').Value = '  -0.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.568'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('E11').Value = '  -3.23%  '
$ws.Range('E12').Value = '  -3.96%  '
$ws.Range('D13').Value = '3.987.03'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.89'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000171'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -8.81%  '
$ws.Range('D17').Value = '63.569.58'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '3.409.27'
$ws.Range('E18').Value = '  -0.16%  '
$ws.Range('E19').Value = '  -4.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.72'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.515'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.50%  '
$ws.Range('E26').Value = '  -4.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.57%  '
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('E31').Value = '  -8.02%  '
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.90'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.50'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.91'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  +8.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('D39').Value = '2.810.97'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.02'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '42.88'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0718'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.40'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.68%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.07%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0303'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '326.12'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('E48').Value = '  +7.37%  '
$ws.Range('E49').Value = '  -4.99%  '
$ws.Range('E50').Value = '  -5.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.57%  '
